$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
# Row 43
$ws.Cells.Item(43, 8).Value = 2205.2666
$ws.Cells.Item(43, 10).Value = 2205.2666
$ws.Cells.Item(43, 12).Value = 2205.2666
$ws.Cells.Item(43, 14).Value = -2343.2666

# Row 132
$ws.Cells.Item(132, 8).Value = 820352
$ws.Cells.Item(132, 9).Value = 3679.8958
$ws.Cells.Item(132, 10).Value = 4087040.5
$ws.Cells.Item(132, 11).Value = 11039.6874
$ws.Cells.Item(132, 12).Value = 12261121.5
$ws.Cells.Item(132, 13).Value = -8509.687399999999
$ws.Cells.Item(132, 14).Value = -12266181.5

# Row 138
$ws.Cells.Item(138, 8).Value = 4004137
$ws.Cells.Item(138, 10).Value = 4765949
$ws.Cells.Item(138, 12).Value = 14297847
$ws.Cells.Item(138, 14).Value = -14308127

$ws = $wb.Worksheets.Item(2)
# Row 2
$ws.Cells.Item(2, 8).Value = 2325.7666
$ws.Cells.Item(2, 9).Value = 2071.1667
$ws.Cells.Item(2, 10).Value = 2707.6667
$ws.Cells.Item(2, 11).Value = 2071.1667
$ws.Cells.Item(2, 12).Value = 2707.6667
$ws.Cells.Item(2, 13).Value = -1958.1667
$ws.Cells.Item(2, 14).Value = -2933.6667

# Row 45
$ws.Cells.Item(45, 8).Value = 1890.6666
$ws.Cells.Item(45, 9).Value = 1836
$ws.Cells.Item(45, 11).Value = 1836
$ws.Cells.Item(45, 13).Value = -1459

# Row 74
$ws.Cells.Item(74, 8).Value = 6212050
$ws.Cells.Item(74, 9).Value = 8359763
$ws.Cells.Item(74, 10).Value = 75727.64
$ws.Cells.Item(74, 11).Value = 8359763
$ws.Cells.Item(74, 12).Value = 75727.64
$ws.Cells.Item(74, 13).Value = -8358889
$ws.Cells.Item(74, 14).Value = -77475.64

# Row 77
$ws.Cells.Item(77, 8).Value = 6212050
$ws.Cells.Item(77, 9).Value = 8359763
$ws.Cells.Item(77, 10).Value = 75727.64
$ws.Cells.Item(77, 11).Value = 41798815
$ws.Cells.Item(77, 12).Value = 378638.2
$ws.Cells.Item(77, 13).Value = -41794447
$ws.Cells.Item(77, 14).Value = -387374.2

# Row 97
$ws.Cells.Item(97, 8).Value = 2842178
$ws.Cells.Item(97, 9).Value = 3290416.8
$ws.Cells.Item(97, 10).Value = 3333.3333
$ws.Cells.Item(97, 11).Value = 3290416.8
$ws.Cells.Item(97, 12).Value = 3333.3333
$ws.Cells.Item(97, 13).Value = -3289920.8
$ws.Cells.Item(97, 14).Value = -4325.3333

# Row 102
$ws.Cells.Item(102, 8).Value = 28573588
$ws.Cells.Item(102, 9).Value = 28573588
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 28573588
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 14).Value = ""
$ws.Cells.Item(102, 13).Value = -28571966

# Row 110
$ws.Cells.Item(110, 8).Value = 2253.1428
$ws.Cells.Item(110, 9).Value = 595.8
$ws.Cells.Item(110, 10).Value = 3173.889
$ws.Cells.Item(110, 11).Value = 595.8
$ws.Cells.Item(110, 12).Value = 3173.889
$ws.Cells.Item(110, 13).Value = 1449.2
$ws.Cells.Item(110, 14).Value = -7263.889

# Row 116
$ws.Cells.Item(116, 8).Value = 2325.7666
$ws.Cells.Item(116, 9).Value = 2071.1667
$ws.Cells.Item(116, 10).Value = 2707.6667
$ws.Cells.Item(116, 11).Value = 2071.1667
$ws.Cells.Item(116, 12).Value = 2707.6667
$ws.Cells.Item(116, 13).Value = 222.8332999999998
$ws.Cells.Item(116, 14).Value = -7295.6667

# Row 132
$ws.Cells.Item(132, 8).Value = 10248606
$ws.Cells.Item(132, 9).Value = 14316024
$ws.Cells.Item(132, 10).Value = 80058.57000000001
$ws.Cells.Item(132, 11).Value = 42948072
$ws.Cells.Item(132, 12).Value = 240175.71
$ws.Cells.Item(132, 13).Value = -42945542
$ws.Cells.Item(132, 14).Value = -245235.71

$ws = $wb.Worksheets.Item(3)
# Row 3
$ws.Cells.Item(3, 8).Value = 2325.7666
$ws.Cells.Item(3, 9).Value = 2071.1667
$ws.Cells.Item(3, 10).Value = 2707.6667
$ws.Cells.Item(3, 11).Value = 2071.1667
$ws.Cells.Item(3, 12).Value = 2707.6667
$ws.Cells.Item(3, 13).Value = -1957.1667
$ws.Cells.Item(3, 14).Value = -2935.6667

# Row 20
$ws.Cells.Item(20, 8).Value = 1678.3334
$ws.Cells.Item(20, 9).Value = 1126.6666
$ws.Cells.Item(20, 10).Value = 1954.1666
$ws.Cells.Item(20, 11).Value = 1126.6666
$ws.Cells.Item(20, 12).Value = 1954.1666
$ws.Cells.Item(20, 13).Value = -879.6666
$ws.Cells.Item(20, 14).Value = -2448.1666

# Row 29
$ws.Cells.Item(29, 8).Value = 886
$ws.Cells.Item(29, 9).Value = 886
$ws.Cells.Item(29, 11).Value = 886
$ws.Cells.Item(29, 13).Value = -597

# Row 80
$ws.Cells.Item(80, 8).Value = 306.4643
$ws.Cells.Item(80, 9).Value = 220
$ws.Cells.Item(80, 10).Value = 325.26086
$ws.Cells.Item(80, 11).Value = 220
$ws.Cells.Item(80, 12).Value = 325.26086
$ws.Cells.Item(80, 13).Value = 778
$ws.Cells.Item(80, 14).Value = -2321.26086

# Row 83
$ws.Cells.Item(83, 8).Value = 306.4643
$ws.Cells.Item(83, 9).Value = 220
$ws.Cells.Item(83, 10).Value = 325.26086
$ws.Cells.Item(83, 11).Value = 1100
$ws.Cells.Item(83, 12).Value = 1626.3043
$ws.Cells.Item(83, 13).Value = 3892
$ws.Cells.Item(83, 14).Value = -11610.3043

# Row 94
$ws.Cells.Item(94, 8).Value = 536.8125
$ws.Cells.Item(94, 9).Value = 521.46155
$ws.Cells.Item(94, 10).Value = 603.3333
$ws.Cells.Item(94, 11).Value = 521.46155
$ws.Cells.Item(94, 12).Value = 603.3333
$ws.Cells.Item(94, 13).Value = -70.46154999999999
$ws.Cells.Item(94, 14).Value = -1505.3333

# Row 99
$ws.Cells.Item(99, 8).Value = 1143.6
$ws.Cells.Item(99, 9).Value = 945.7143
$ws.Cells.Item(99, 10).Value = 1605.3334
$ws.Cells.Item(99, 11).Value = 945.7143
$ws.Cells.Item(99, 12).Value = 1605.3334
$ws.Cells.Item(99, 13).Value = 552.2857
$ws.Cells.Item(99, 14).Value = -4601.3334

# Row 107
$ws.Cells.Item(107, 8).Value = 2266.3333
$ws.Cells.Item(107, 9).Value = 2266.3333
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 2266.3333
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 14).Value = ""
$ws.Cells.Item(107, 13).Value = -346.3332999999998

# Row 134
$ws.Cells.Item(134, 8).Value = 10755850
$ws.Cells.Item(134, 9).Value = 2980.95
$ws.Cells.Item(134, 11).Value = 8942.849999999999
$ws.Cells.Item(134, 13).Value = -6407.849999999999

$ws = $wb.Worksheets.Item(4)
# Row 16
$ws.Cells.Item(16, 8).Value = 2166.6667
$ws.Cells.Item(16, 9).Value = 2500
$ws.Cells.Item(16, 10).Value = 2000
$ws.Cells.Item(16, 11).Value = 2500
$ws.Cells.Item(16, 12).Value = 2000
$ws.Cells.Item(16, 13).Value = -2213
$ws.Cells.Item(16, 14).Value = -2574

# Row 58
$ws.Cells.Item(58, 8).Value = 55556750
$ws.Cells.Item(58, 9).Value = 55556750
$ws.Cells.Item(58, 11).Value = 55556750
$ws.Cells.Item(58, 13).Value = -55556547

# Row 93
$ws.Cells.Item(93, 8).Value = 35000
$ws.Cells.Item(93, 9).Value = 35000
$ws.Cells.Item(93, 11).Value = 35000
$ws.Cells.Item(93, 13).Value = -33128

# Row 105
$ws.Cells.Item(105, 8).Value = 667.375
$ws.Cells.Item(105, 9).Value = 486.16666
$ws.Cells.Item(105, 10).Value = 1211
$ws.Cells.Item(105, 11).Value = 486.16666
$ws.Cells.Item(105, 12).Value = 1211
$ws.Cells.Item(105, 13).Value = 1260.83334
$ws.Cells.Item(105, 14).Value = -4705

# Row 113
$ws.Cells.Item(113, 8).Value = 2166.6667
$ws.Cells.Item(113, 9).Value = 2500
$ws.Cells.Item(113, 10).Value = 2000
$ws.Cells.Item(113, 11).Value = 2500
$ws.Cells.Item(113, 12).Value = 2000
$ws.Cells.Item(113, 13).Value = -330
$ws.Cells.Item(113, 14).Value = -6340

# Row 134
$ws.Cells.Item(134, 8).Value = 84939.69500000001
$ws.Cells.Item(134, 9).Value = 1968.4445
$ws.Cells.Item(134, 10).Value = 271625
$ws.Cells.Item(134, 11).Value = 5905.333500000001
$ws.Cells.Item(134, 12).Value = 814875
$ws.Cells.Item(134, 13).Value = -3370.333500000001
$ws.Cells.Item(134, 14).Value = -819945

# Row 136
$ws.Cells.Item(136, 8).Value = 55556750
$ws.Cells.Item(136, 9).Value = 55556750
$ws.Cells.Item(136, 11).Value = 166670250
$ws.Cells.Item(136, 13).Value = -166667700

$ws = $wb.Worksheets.Item(5)
# Row 9
$ws.Cells.Item(9, 8).Value = 51000
$ws.Cells.Item(9, 10).Value = 2000
$ws.Cells.Item(9, 12).Value = 6000
$ws.Cells.Item(9, 14).Value = -6448

# Row 113
$ws.Cells.Item(113, 8).Value = 723.3158
$ws.Cells.Item(113, 9).Value = 593.3333
$ws.Cells.Item(113, 10).Value = 840.3
$ws.Cells.Item(113, 11).Value = 1779.9999
$ws.Cells.Item(113, 12).Value = 2520.9
$ws.Cells.Item(113, 13).Value = 390.0001
$ws.Cells.Item(113, 14).Value = -6860.9

$ws = $wb.Worksheets.Item(6)
# Row 113
$ws.Cells.Item(113, 8).Value = 2599.8235
$ws.Cells.Item(113, 9).Value = 1614.625
$ws.Cells.Item(113, 10).Value = 3475.5557
$ws.Cells.Item(113, 11).Value = 1614.625
$ws.Cells.Item(113, 12).Value = 3475.5557
$ws.Cells.Item(113, 13).Value = 555.375
$ws.Cells.Item(113, 14).Value = -7815.5557

# Row 126
$ws.Cells.Item(126, 8).Value = 1907.4445
$ws.Cells.Item(126, 9).Value = 1714
$ws.Cells.Item(126, 10).Value = 2149.25
$ws.Cells.Item(126, 11).Value = 5142
$ws.Cells.Item(126, 12).Value = 6447.75
$ws.Cells.Item(126, 13).Value = -2672
$ws.Cells.Item(126, 14).Value = -11387.75

# Row 132
$ws.Cells.Item(132, 8).Value = 337749.66
$ws.Cells.Item(132, 9).Value = 506000
$ws.Cells.Item(132, 10).Value = 253624.5
$ws.Cells.Item(132, 11).Value = 1518000
$ws.Cells.Item(132, 12).Value = 760873.5
$ws.Cells.Item(132, 13).Value = -1515470
$ws.Cells.Item(132, 14).Value = -765933.5

$ws = $wb.Worksheets.Item(7)
# Row 19
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 14).Value = ""

# Row 46
$ws.Cells.Item(46, 8).Value = 3788449
$ws.Cells.Item(46, 9).Value = 7576298
$ws.Cells.Item(46, 10).Value = 600
$ws.Cells.Item(46, 11).Value = 7576298
$ws.Cells.Item(46, 12).Value = 600
$ws.Cells.Item(46, 13).Value = -7576110
$ws.Cells.Item(46, 14).Value = -976

# Row 136
$ws.Cells.Item(136, 8).Value = 54836.5
$ws.Cells.Item(136, 9).Value = 35982.97
$ws.Cells.Item(136, 10).Value = 130250.625
$ws.Cells.Item(136, 11).Value = 107948.91
$ws.Cells.Item(136, 12).Value = 390751.875
$ws.Cells.Item(136, 13).Value = -105398.91
$ws.Cells.Item(136, 14).Value = -395851.875

$ws = $wb.Worksheets.Item(8)
# Row 81
$ws.Cells.Item(81, 8).Value = 2287.35
$ws.Cells.Item(81, 9).Value = 690
$ws.Cells.Item(81, 10).Value = 2819.8
$ws.Cells.Item(81, 11).Value = 1380
$ws.Cells.Item(81, 12).Value = 5639.6
$ws.Cells.Item(81, 13).Value = -319
$ws.Cells.Item(81, 14).Value = -7761.6

# Row 84
$ws.Cells.Item(84, 8).Value = 2287.35
$ws.Cells.Item(84, 9).Value = 690
$ws.Cells.Item(84, 10).Value = 2819.8
$ws.Cells.Item(84, 11).Value = 6900
$ws.Cells.Item(84, 12).Value = 28198
$ws.Cells.Item(84, 13).Value = -1596
$ws.Cells.Item(84, 14).Value = -38806

# Row 132
$ws.Cells.Item(132, 8).Value = 50423.83
$ws.Cells.Item(132, 9).Value = 42849.082
$ws.Cells.Item(132, 11).Value = 128547.246
$ws.Cells.Item(132, 13).Value = -126017.246

# Row 136
$ws.Cells.Item(136, 8).Value = 43572.49
$ws.Cells.Item(136, 9).Value = 31805.848
$ws.Cells.Item(136, 10).Value = 67841.19
$ws.Cells.Item(136, 11).Value = 95417.54400000001
$ws.Cells.Item(136, 12).Value = 203523.57
$ws.Cells.Item(136, 13).Value = -92867.54400000001
$ws.Cells.Item(136, 14).Value = -208623.57
